$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$questions = @{}
$questions[2] = 'کدام عضله گردن دارای دو عصب دهی است؟'
$questions[3] = 'معدنی شدن دندان ها از کجا شروع می شود؟ comedk 07'
$questions[4] = 'یک کارگر آنگنوادی برای چه مدت آموزش می‌بیند؟'
$questions[5] = 'مأموریت ملی سلامت روستایی توسط نخست وزیر در چه سالی اعلام شد؟ سپتامبر 2007'
$questions[6] = 'یک مرد 78 ساله از بیاشتهایی، خستگی و به طور کلی احساس ناخوشی شکایت دارد. در معاینه، طحال بزرگ، رنگپریدگی مشاهده میشود و بقیه معاینه طبیعی است. شمارش خون او غیرطبیعی است: wbc برابر با 46000 در میلیلیتر با افزایش گرانولوسیتهای بالغ، هموگلوبین 9.0 گرم در دسیلیتر و پلاکت‌ها 450,000 در میلیلیتر. بررسی‌های بیشتر، رونوشت‌های bcr-abl را نشان می‌دهد. همچنین جهش t315i مشاهده شده است. داروی انتخابی برای این وضعیت کدام است؟'
$questions[7] = 'در مقایسه با هپارین، انوکساپارین:'
$questions[8] = 'در سونوگرافی کبد، الگوی "آسمان ستاره‌ای" دیده می‌شود. این یافته ویژگی کدام یک از موارد زیر است؟'
$questions[9] = 'نمره‌دهی پیرانی برای ctev شامل همه موارد زیر به جز کدام است؟'
$questions[10] = 'کدام نوع از بیماری ذخیره‌ای گلیکوژن بر عضلات تأثیر نمی‌گذارد؟'
$questions[11] = 'زمانی که عامل بیماری‌زا وجود دارد اما انتقالی صورت نمی‌گیرد، به چه چیزی معروف است؟'
$questions[12] = 'همه موارد زیر از ویژگی‌های آلوئولیت (بیماری بینابینی ریه) هستند، به جز -'
$questions[13] = 'سیس آتراکوریوم به دلیل کدام مزیت نسبت به آتراکوریوم ترجیح داده می‌شود؟'
$questions[14] = 'شایع‌ترین نوع اتصال غیرطبیعی کامل وریدهای ریوی کدام است؟'
$questions[15] = 'گانگرین فورنیر در کدام ناحیه مشاهده می‌شود؟ سپتامبر 2008'
$questions[16] = 'یک زن در شب احساس درد و حس خزیدن حشره روی پاهای خود دارد که با تکان دادن پاها تسکین مییابد. کدام یک از گزینههای زیر داروی انتخابی برای این شرایط است؟'
$questions[17] = 'کدام یک از عبارات زیر در مورد siadh صحیح است؟'
$questions[18] = 'در یک کارآزمایی کنترل شده برای مقایسه دو روش درمانی، هدف اصلی تصادفی‌سازی این است که اطمینان حاصل شود؟'
$questions[19] = 'کدام یک از موارد زیر در مورد فرورفتگی صورت (facial recess) صحیح است؟'
$questions[20] = 'تفاوت بین عمل ویروس و باکتری در چیست؟'
$questions[21] = 'بهترین درمان برای sbp در شرایط سیروز و آسیت، تزریق سفوتاکسیم است. بهترین جایگزین خوراکی برای درمان آن کدام است؟'
$questions[22] = 'رابر دام توصیه شده برای ایجاد رتراکشن مؤثر بافت لثه کدام است؟'
$questions[23] = 'زخم شدن لکه‌های پییر در عفونت _________ رخ می‌دهد'
$questions[24] = 'بیمار با از دست دادن حس در سطوح مجاور انگشت شست و دوم پا و اختلال در دورسی فلکشن پا مراجعه می‌کند. این علائم احتمالاً نشان‌دهنده آسیب به کدام یک از اعصاب زیر است؟'
$questions[25] = 'همه موارد زیر از ویژگی‌های سندرم کورساکف هستند به جز'
$questions[26] = 'کدام یک از موارد زیر از عوارض لنفادم مزمن اندام محسوب نمی‌شود؟'
$questions[27] = 'در کدام یک از موارد زیر، فشار داخل چشمی بسیار بالا و التهاب حداقل است؟'
$questions[28] = 'استئوکلاستوما معمولاً در کدام قسمت رخ می‌دهد؟'
$questions[29] = 'اثرات سیستمیک سندرم آپنه انسدادی خواب (osas) شامل همه موارد زیر به جز کدام است؟'
$questions[30] = 'کدام یک از داروهای زیر جزو "درمان سه گانه" سرکوب کننده سیستم ایمنی برای بیماران پس از پیوند کلیه نیست؟'
$questions[31] = 'بازدارنده‌های غیررقابتی'
$questions[32] = 'dash چیست؟'
$questions[33] = 'صدای بروئی در غده تیروئید در کجا شنیده می‌شود؟'
$questions[34] = 'در آنژیوگرافی شریان شبکیه، رنگ از طریق کدام محل تزریق می‌شود؟'
$questions[35] = 'علل طولانی شدن مرحله اول زایمان همه موارد زیر هستند به جز: مارس 2009'
$questions[36] = 'رنگ سیلندر اتیلن چیست؟'
$questions[37] = 'غضروف الاستیک در کدام قسمت دیده می‌شود؟'
$questions[38] = 'استفاده طولانی مدت از پتیدین اجتناب می‌شود زیرا متابولیت پتیدین با کدام یک از موارد زیر مرتبط است؟'
$questions[39] = 'واحد متابولیک پایه استخوان چیست؟'
$questions[40] = 'در کدام یک از مراحل زیر atp آزاد می‌شود؟'
$questions[41] = 'شایع‌ترین عارضه ect اصلاح‌شده چیست؟'
$questions[42] = 'شکستگی انفجاری مدار چشم، اغلب منجر به شکستگی کدام ناحیه می‌شود؟'
$questions[43] = 'نام شریان مشخص شده چیست؟'
$questions[44] = 'معده مخملی قرمز در مسمومیت با کدام ماده دیده می‌شود؟'
$questions[45] = 'اشک تمساح به چه علتی رخ می‌دهد؟'
$questions[46] = 'کدام یک از موارد زیر در مورد pnh نادرست است؟'
$questions[47] = 'عدد دیبوکین به چه چیزی اشاره دارد؟'
$questions[48] = 'یک زن ۳۶ ساله متوجه عدم قاعدگی در طول ۴ ماه گذشته شده است. تست بارداری منفی است. سطح هورمون لوتئینیزه کننده و هورمون محرک فولیکول در سرم افزایش یافته و سطح استرادیول سرم پایین است. این یافته‌ها نشان‌دهنده کدام مورد است؟'
$questions[49] = 'نیترات ها مصرف اکسیژن میوکارد را از طریق تمام مکانیسم های زیر کاهش می دهند به جز -'
$questions[50] = 'کدام یک از شرایط زیر خطر افزایش یافته‌ای برای کارسینوم مهاجم پستان ندارد؟'
$questions[51] = 'کدام یک از اعصاب جمجمه‌ای زیر را می‌توان در سی‌تی اسکن ساده مشاهده کرد؟'
$questions[52] = 'موارد منع معاینه مقعدی (dre)'
$questions[53] = 'مهم‌ترین تظاهر فیستول حالب-واژن چیست؟ سپتامبر 2009'
$questions[54] = 'در رادیوگرافی، شکم بدون گاز در کدام یک از موارد زیر دیده می‌شود؟'
$questions[55] = 'کدام یک از موارد زیر یک ریبوزیم است؟'
$questions[56] = 'بیمار و روان درمانگر، هر دو به صورت فعال در کدام یک مشارکت می‌کنند؟'
$questions[57] = 'قطعی‌ترین روش برای تشخیص آمبولی ریه کدام است؟'
$questions[58] = 'کد رنگی کیسه‌های بیمارستانی برای دفع ضایعات آناتومیک انسانی مانند اعضای بدن چیست؟'
$questions[59] = 'یک زن 40 ساله با تنگی نفس پیشرونده از یک سال قبل. همچنین کریپیتاسیون دوطرفه اند-انسپیراتوری در قاعده ریه دارد. گزینه‌های صحیح عبارتند از:  
a) ممکن است با بیماری بافت همبند همراه باشد  
b) حجم باقی‌مانده افزایش یافته است  
c) ظرفیت کل ریوی افزایش یافته است  
d) hrct یک تست تشخیصی مفید است  
e) نسبت fev1/fvc افزایش یافته است'
$questions[60] = 'نکروز کازئوز در کدام یک از موارد زیر دیده می‌شود؟'
$questions[61] = 'مخاط واژن با چه نوع اپیتلیومی پوشیده شده است؟'
$questions[62] = 'یک بیمار 30 ساله g3p2 در هفته 20 بارداری به کلینیک مراقبت های دوران بارداری مراجعه می کند. او در شرح حال بیان می کند که اولین نوزادش با وزن 4.6 کیلوگرم به روش سزارین به دنیا آمده و دومین نوزاد با وزن 4.8 کیلوگرم به روش سزارین متولد شده است. متخصص زنان به دیابت بارداری مشکوک شده و تست تحمل گلوکز (gct) را درخواست می کند. سطح قند خون پس از مصرف 50 گرم گلوکز خوراکی 206 میلی گرم در دسی لیتر است و بنابراین بیمار به عنوان یک مورد دیابت بارداری تأیید می شود. همه موارد زیر از عوارض شناخته شده این وضعیت هستند به جز:'
$questions[63] = 'اولین نشانه کمبود ویتامین a چیست؟'
$questions[64] = 'سندرم گیتلمان در تمام موارد زیر به جز کدام مورد با سندرم باتر تفاوت دارد؟'
$questions[65] = 'کدام گزینه در مورد فک پایین (ماندیبل) صحیح نیست؟'
$questions[66] = 'ترمیم فاترگیل همچنین به عنوان چه چیزی شناخته می‌شود؟'
$questions[67] = 'مردی با اورتریت و ترشح مجرای ادرار به کلینیک بیماری‌های مقاربتی مراجعه می‌کند. رنگ‌آمیزی گرم نشان‌دهنده تعداد زیادی سلول چرکی اما عدم وجود میکروارگانیسم است. کشت در محیط‌های معمول آزمایشگاهی منفی است. محتمل‌ترین عامل کدام است؟'
$questions[68] = 'سفتی نعشی ابتدا در کدام قسمت مشاهده می‌شود؟'
$questions[69] = 'اولین بار جفری در یک پرونده جنایی از اثر انگشت dna برای تشخیص چه موردی استفاده کرد؟ mahe 11'
$questions[70] = 'کم آبی هیپرناترمی با کدام مورد(ها) مشخص می‌شود؟  
الف) سدیم سرمی بیشتر از 150 میلی‌مول بر لیتر  
ب) علائم کم آبی حداقل هستند  
ج) کاهش حجم مایع خارج سلولی  
د) نیاز به اصلاح سریع دارد  
ه) انتقال آب از مایع خارج سلولی به داخل سلولی'
$questions[71] = 'در دندانپزشکی چهار دستی، پاهای دندانپزشک در طول هر پروسه کجا قرار می‌گیرد؟'
$questions[72] = 'در nrhm، کارگران asha از همان روستا استخدام می‌شوند. کدام یک از اصول مراقبت‌های اولیه سلامت در اینجا رعایت شده است؟'
$questions[73] = 'سمیت بلئومایسین چه نوع سلول هایی را تحت تأثیر قرار می دهد؟'
$questions[74] = 'alteplase با چه مکانیسمی عمل می‌کند؟'
$questions[75] = 'یک نوزاد در زمان مراجعه به اورژانس، تب بالا و ناراحتی تنفسی داشت. نمونه خلط، کوکسی‌های گرم مثبت با فعالیت آلفا همولیتیک نشان داد. عامل بیماری‌زای مشکوک به کدام یک از عوامل زیر حساس است؟'
$questions[76] = 'همه موارد زیر با درمان بیکربنات سدیم مشاهده می‌شوند به جز:'
$questions[77] = 'بیمار با سابقه ضربه غیرنافذ به قفسه سینه و شکم مراجعه کرده است. در بررسی اولتراسوند با حالت m یافته زیر مشاهده می‌شود. احتمالاً تشخیص چیست؟'
$questions[78] = 'تکنیک کوله سیستکتومی لاپاراسکوپی برای اولین بار توسط چه کسی توصیف شد؟'
$questions[79] = 'در کدام یک از بیماری‌های ذخیره‌سازی گلیکوژن، عضلات درگیر نیستند؟'
$questions[80] = 'کوردوم نازوفارنکس از کجا منشأ می‌گیرد؟'
$questions[81] = 'استرپتوکوک پنومونیا کدام نوع همولیز را تولید می‌کند؟'
$questions[82] = 'یک نوزاد 5 ماهه با سابقه استفراغ، تحریک‌پذیری و زردی به اورژانس مراجعه می‌کند. نوزاد قبلاً سالم بود. سابقه معرفی آبمیوه به رژیم غذایی یک هفته قبل وجود دارد. سونوگرافی هپاتومگالی را نشان داد و آزمایش‌های عملکرد کبد غیرطبیعی بودند. نوزاد به احتمال زیاد کمبود کدام یک از آنزیم‌های زیر را دارد؟'
$questions[83] = 'کدام یک از تومورها منشأ ویروسی دارند؟'
$questions[84] = 'جدی‌ترین عارضه شکستگی استخوان بلند چیست؟ سپتامبر 2005'
$questions[85] = 'بیشترین غلظت پروستاگلاندین در کدام مورد دیده می‌شود؟'
$questions[86] = 'در تصویر رادیوگرافی جمجمه، کدام یک از ضایعات مغزی مشاهده می‌شود؟'
$questions[87] = 'فضای پری پلاسمی در کدام یک دیده می‌شود؟'
$questions[88] = 'عوامل خطر سندرم تخمدان پلی کیستیک قبل از بلوغ همه موارد زیر هستند به جز'
$questions[89] = 'تمامی عبارات زیر در مورد کروموزوم فیلادلفیا در cml صحیح هستند، به جز:'
$questions[90] = 'فاگوسیتوز کریپتوکوکوس نئوفورمانس توسط چه چیزی مهار می‌شود؟'
$questions[91] = 'مجاری بلینی در کدام اندام یافت می‌شوند؟'
$questions[92] = 'همه موارد زیر در مورد کانکروم اوریس صحیح است به جز -'
$questions[93] = 'کدام بیماری قلبی-عروقی در دوران بارداری بالاترین میزان مرگ و میر مادران را دارد؟'
$questions[94] = 'یک مرد 65 ساله با فلج تارهای صوتی مراجعه کرده است. او از کودکی عادت به جویدن پان دارد. در معاینه، عملکرد abduction تارهای صوتی او مختل شده است. عصب درگیر کدام است؟'
$questions[95] = 'درگیری علامت‌دار سیستم عصبی مرکزی در نوزادان بیشتر در کدام گروه از عفونت‌های مادرزادی داخل رحمی دیده می‌شود؟'
$questions[96] = 'ارزیابی حجم مایع خارج سلولی (ecf) را می‌توان به دقیق‌ترین شکل با کدام ماده انجام داد؟'
$questions[97] = 'کدام یک از غلات را می‌توان با اطمینان در بیماری اسپروی سلیاک استفاده کرد؟  
الف) ذرت  
ب) چاودار  
ج) جو  
د) برنج'
$questions[98] = 'خون‌سازی ابتدا در کدام قسمت شروع می‌شود؟'
$questions[99] = 'بیمار یک ماه پس از بیوپسی خوش‌خیم پستان راست، با یک طناب زیرپوستی جانبی که دقیقاً زیر پوست حس می‌شود و باعث درد می‌شود، مراجعه می‌کند. علت این وضعیت چیست؟'
$questions[100] = 'پسوریازیس که معمولاً پس از عفونت استرپتوکوکی β-همولیتیک مشاهده می‌شود کدام است؟'
$questions[101] = 'تنگی عروق ریوی محیطی در تمام موارد زیر دیده می‌شود به جز'
$questions[102] = 'یک نوزاد در حالت شوک در اورژانس حاضر می‌شود. کدام مایع ایده‌آل برای تجویز به نوزاد است؟'
$questions[103] = 'کدام پیوند شیمیایی در کراتین پوست و ناخن وجود دارد که باعث تفاوت در قوام آنها می‌شود؟'
$questions[104] = 'اثر توبوکورارین به راحتی توسط کدام یک از موارد زیر برعکس می‌شود؟'
$questions[105] = 'هر یک از عبارات زیر را در مورد مایع داخل سلولی که صحیح است، بررسی کنید:  
1. بیش از 50٪ از آب بدن را شامل می‌شود  
2. فشار اسمزی بالاتری نسبت به مایع خارج سلولی دارد  
3. غلظت آنیون‌های آلی بیشتری نسبت به مایع خارج سلولی دارد  
کدام یک از گزینه‌های زیر بهترین نتیجه‌گیری را نشان می‌دهد؟'
$questions[106] = 'بدخیمی پنهان تیروئید با متاستازهای گره‌ای در کدام مورد دیده می‌شود؟'
$questions[107] = 'دندان‌های بدون ریشه به چه شرایطی اشاره دارند؟'
$questions[108] = 'کدام یک از موارد زیر جزء تشخیص‌های افتراقی ندول منفرد ریوی محسوب نمی‌شود؟'
$questions[109] = 'شستشوی معده در مسمومیت با کدام یک قابل انجام است؟ wb 10; odisha 11'
$questions[110] = 'بیمار با زخم بدون درد روی آلت تناسلی خود مراجعه می‌کند. تست واسرمن انجام شده و مثبت است. واکنش واسرمن نمونه‌ای از ________ است.'
$questions[111] = 'تعریف نرخ مرگ و میر مادران شامل همه موارد زیر به جز کدام است؟'
$questions[112] = 'کمبود کدام ویتامین باعث گلوسیت و کیلوزیس می‌شود؟'
$questions[113] = 'شایع‌ترین سرطان مقعد کدام است؟'
$questions[114] = 'شایع‌ترین داروی ایجادکننده انتروکولیت سودوممبران چیست؟'
$questions[115] = 'کف پای نوزاد تازه متولد شده که در زیر نشان داده شده است، نشانگر کدام مورد است؟'
$questions[116] = 'کدام یک از بیهوش کننده های عمومی زیر اثر شل کنندگی عضلانی ضعیفی دارد؟'
$questions[117] = 'کدام یک از موارد زیر واکسن زنده ضعیف شده است؟'
$questions[118] = 'کدام مورد در فلج مغزی یافت نمی‌شود؟'
$questions[119] = 'یک زن ۵۰ ساله توسط پزشک خانواده اش ارجاع داده شده است زیرا از اختلال افسردگی رنج می‌برد. در طول مصاحبه، او در بیان احساسات خود مشکل دارد. این پدیده به بهترین شکل چگونه توصیف می‌شود؟'
$questions[120] = 'آدنوم پلئومورفیک از کدام یک منشأ می‌گیرد؟'
$questions[121] = 'گیرنده هورمون تیروئید که عمدتاً در سیستم عصبی مرکزی/مغز بیان می‌شود -'
$questions[122] = 'کودکی با تب به مدت ۲ روز، تغییر سطح هوشیاری و راش های پورپوریک مراجعه کرده است. فشار خون او ۹۰/۶۰ mmhg است. درمان انتخابی چیست؟'
$questions[123] = 'یک مرد 24 ساله به مدت سه ماه گذشته چندین لکه کوچک هیپوپیگمانته روی قفسه سینه و پشت داشته است. این لکه‌ها گرد بوده، حول فولیکول‌ها قرار گرفته‌اند و بسیاری از آن‌ها به هم پیوسته و صفحات بزرگی تشکیل داده‌اند. سطح لکه‌ها پوسته‌ریزی ظریفی دارد. او یک سال پیش نیز ضایعات مشابهی داشته که با درمان بهبود یافته‌اند. مناسب‌ترین روش بررسی برای تأیید تشخیص چیست؟'
$questions[124] = 'ژن مسئول انتقال اسید فولیک در کدام کروموزوم قرار دارد؟'
$questions[125] = 'لکه مک کالوم نشانه تشخیصی کدام بیماری است؟'
$questions[126] = '"گردن گاوی" در موارد شدید کدام یک از موارد زیر دیده می‌شود؟'
$questions[127] = 'کدام یک از موارد زیر شایع‌ترین ناهنجاری عروقی کلیوی است؟'
$questions[128] = 'یک زن ۵۲ ساله برای معاینه به پزشک مراجعه می‌کند. او در حال بهبودی از شکستگی مچ دست پس از یک زمین‌خوردن است. اسکن جذب‌سنجی اشعه ایکس با انرژی دوگانه (dexa) از لگن نشان‌دهنده پوکی استخوان بوده است. او در سن ۵۰ سالگی یائسه شده و به دلیل سابقه خانوادگی قوی سرطان پستان، درمان جایگزینی هورمون را شروع نکرده است. اکنون او از شکستگی لگن در آینده می‌ترسد و می‌خواهد یک رژیم پیشگیری از کاهش استخوان را شروع کند. کدام یک از عوامل دارویی زیر برای این بیمار مناسب‌تر است؟'
$questions[129] = 'استانداردسازی مستقیم برای مقایسه نرخ‌های مرگ و میر بین دو کشور استفاده می‌شود. این کار به دلیل تفاوت در کدام یک از پارامترهای زیر انجام می‌شود؟'
$questions[130] = 'قطر دو گیجگاهی جنین چقدر است؟'
$questions[131] = 'تمام داروهای زیر پس بار را کاهش می‌دهند، به جز'
$questions[132] = 'طول مری در بزرگسالان چقدر است؟'
$questions[133] = 'چه کسی واکنش زنجیره ای پلیمراز را کشف کرد؟'
$questions[134] = 'اسهال ترشحی در کدام یک از موارد زیر دیده نمی‌شود؟'
$questions[135] = 'کدام یک از دردهای زیر به گوش ارجاع داده نمی شود؟'
$questions[136] = 'بیماری پاژه پستان با کدام یک از موارد زیر مرتبط است؟'
$questions[137] = 'بازداری رنشاو'
$questions[138] = 'بی حسی اپیدورال برای کدام مورد مناسب است؟'
$questions[139] = 'رنگدانه‌های پوستی شبیه به «ردپای راه‌آهن» در کدام مورد دیده می‌شود؟'
$questions[140] = 'در مورد هیبریدوما همه موارد زیر صحیح است به جز؟'
$questions[141] = 'کدام یک از موارد زیر فلور غالب دهان است؟'
$questions[142] = 'چوب سقط جنین که در سقط جنین جنایی استفاده می‌شود، با چه مکانیزمی باعث سقط جنین می‌شود؟'
$questions[143] = 'عدسی شامل قدیمی ترین سلول ها در کدام قسمت است؟'
$questions[144] = 'کسر تخلیه ای حدوداً چقدر است؟'
$questions[145] = 'یک دندان آسیاب اول فک بالا در دندان های شیری چند ریشه دارد؟'
$questions[146] = 'یک زن ۳۸ ساله با شکایت از کاهش شنوایی در گوش راست طی دو سال گذشته مراجعه کرده است. در تست با چنگال تنظیم ۵۱۲ هرتز، تست رینه بدون ماسک در گوش راست منفی و در گوش چپ مثبت است. در تست وبر، صدا در گوش چپ بلندتر воспринима می‌شود. محتمل‌ترین مشکل بیمار چیست؟'
$questions[147] = 'صفات بی‌احساس و بی‌عاطفه (cu) که شامل بی‌توجهی به دیگران، کمبود همدلی و عموماً عاطفه ناکافی است، به عنوان مشخصه‌ای برای کدام اختلال در dsm v در نظر گرفته می‌شود؟'
$questions[148] = 'ضریب هوشی بین 50 تا 70 نشان دهنده چه چیزی است؟'
$questions[149] = 'درمان انتخابی برای واکنش نوع دوم جذام (لپرا) کدام است؟'
$questions[150] = 'قلیایی کردن ادرار در مدیریت مسمومیت با عوامل زیر مؤثر است:'
$questions[151] = 'موارد استفاده از نور وود شامل موارد زیر است:'
$questions[152] = 'وظایف شغلی دستیاران سلامت مرد چیست؟'
$questions[153] = 'ویروس‌ها را می‌توان با استفاده از تمام موارد زیر پرورش داد به جز:'
$questions[154] = 'زاویه فوقانی کتف در چه سطحی قرار دارد؟'
$questions[155] = 'استخوان سزاموئید در تاندون کدام یک از عضلات زیر وجود دارد؟'
$questions[156] = 'در ارزیابی درد اندام فوقانی، تمام موارد زیر جزو بررسی‌های تشخیصی مفید هستند، به جز:'
$questions[157] = 'کدام یک از ویروس‌های زیر دارای اجسام درون‌هسته‌ای و درون‌سیتوپلاسمی است؟'
$questions[158] = 'گیرنده های آلدوسترون در همه موارد زیر وجود دارند به جز:'
$questions[159] = 'کدام یک از موارد زیر از طریق مکانیسم بازخورد منفی تنظیم نمی‌شود؟'
$questions[160] = 'ایمنی سلولی از نوع کدام است؟'
$questions[161] = 'برش میدلاین در شکم، طول بخیه مورد نیاز-'
$questions[162] = 'کوتاهی متاکارپال در کدام یک دیده می‌شود؟'
$questions[163] = 'همه موارد زیر از مشتقات تاج عصبی هستند به جز؟'
$questions[164] = 'ضخیم‌ترین عصب جمجمه‌ای کدام است؟'
$questions[165] = 'کدام یک از موارد زیر دلیل ریخته‌گری سیاه است؟'
$questions[166] = 'کدام یک از موارد زیر با اندازه گیری شاخص بریدگی سیاتیک در اسکلت تعیین می شود؟'
$questions[167] = 'تمامی موارد زیر با پرکاری تیروئید مرتبط هستند، به جز:'
$questions[168] = 'پوشش اپیتلیال واژن چیست؟'
$questions[169] = 'اتصال به دنده اول شامل همه موارد زیر به جز کدام است؟'
$questions[170] = 'کدام یک از ویژگی های زیر نشان دهنده وجود بیماری قلبی در بارداری است و در بارداری طبیعی دیده نمی شود؟'
$questions[171] = 'کدام یک از استاتین های زیر را می توان در هر زمان از روز مصرف کرد؟'
$questions[172] = 'شایع‌ترین ضایعه پیش‌سرطانی برای بدخیمی دهانی کدام است؟'
$questions[173] = 'گلوکزوری زمانی رخ می‌دهد که غلظت گلوکز خون وریدی از چه مقداری تجاوز کند؟'
$questions[174] = 'وزن کره چشم یک فرد بالغ چقدر است؟'
$questions[175] = 'نیاز انرژی روزانه اضافی در شش ماه اول دوران شیردهی یک زن چقدر است؟'
$questions[176] = 'در مورد asha کدام گزینه صحیح است؟'
$questions[177] = 'سرطان نازوفارنکس شامل کدام نواحی می‌شود؟'
$questions[178] = 'یک ورزشکار پرش با نیزه در حین پرش با نیزه سقوط کرد و دچار فلج بازو شد. کدام یک از بررسی‌های زیر بهترین پیش‌آگهی بهبودی را ارائه می‌دهد؟'
$questions[179] = 'در dna، آدنین با کدام باز جفت می‌شود؟'
$questions[180] = 'کدام یک از مواد زیر برای پوشش غیرمستقیم پالپ ایده‌آل‌تر است؟'
$questions[181] = 'کدام یک از یافته های شبکیه ای زیر در معاینه یک بیمار با فشار خون بدخیم (malignant hypertension) مشخصه است؟'
$questions[182] = 'در هماچوری با منشأ گلومرولی، ادرار با وجود همه موارد زیر به جز کدام یک مشخص می‌شود؟'
$questions[183] = 'یک مرد ۸۲ ساله به دلیل خونریزی شدید از زخم معده دچار شوک هیپوولمیک شده و فوت می‌کند. در کالبدشکافی، نکروز سانتری لوبولار در کبد مشاهده می‌شود. در مقایسه با هپاتوسیت‌های سالم، سلول‌های نکروتیک غلظت داخل سلولی بالاتری از کدام مورد زیر را دارند؟'
$questions[184] = 'شایع‌ترین علت تنگی مجرای ادرار در یک فرد جوان چیست؟'
$questions[185] = 'نادرست در مورد اسپوروتریکس چیست؟'
$questions[186] = 'برای پیوند، قرنیه در چه محیطی نگهداری می‌شود؟'
$questions[187] = 'در مدیریت بیماران مبتلا به عفونت وینسنت، عوامل درمانی ترجیحی کدامند؟'
$questions[188] = 'اسیدهای مونوپروتیک کدامند؟'
$questions[189] = 'بیمار با سنگ‌های کلیوی مکرر مراجعه می‌کند. بررسی میکروسکوپی نمونه ادرار در زیر نشان داده شده است. کدام یک از موارد زیر در ادرار این بیمار دیده نمی‌شود؟'
$questions[190] = 'رنگ زعفرانی مکونیوم در کدام مورد دیده می‌شود؟'
$questions[191] = 'سارکوم کاپوسی در کدام مورد مشاهده می‌شود؟'
$questions[192] = 'برش داخلی مایل فلپ پریودنتال:'
$questions[193] = 'سطح مجاز فلوراید در آب آشامیدنی چقدر است؟'
$questions[194] = 'یک دارو با دوز بارگذاری 20 میلی‌گرم به غلظت پلاسمایی 0.5 میلی‌گرم در لیتر می‌رسد. اگر حجم توزیع ظاهری 40 لیتر باشد، زیست‌دسترسی دارو را محاسبه کنید.'
$questions[195] = 'همه موارد زیر در هرپس زوستر چشمی رخ می‌دهد به جز'
$questions[196] = 'کدام یک از عوامل غیر دپولاریزان بلوک عصبی عضلانی کوتاه ترین اثر را دارد؟'
$questions[197] = 'دیاتوم‌ها در کدام مورد دیده می‌شوند؟'
$questions[198] = 'فالانگا چیست؟'
$questions[199] = 'کمبود تیامین - علل؟'
$questions[200] = 'کدام یک از گزینه‌های زیر داروی خط اول برای القای فولیکول در ناباروری ناشی از سندرم تخمدان پلی‌کیستیک (pcos) است؟'
$questions[201] = 'هنگامی که ارتباط بین دو متغیر توسط یک متغیر سوم به دلیل ارتباط غیرمستقیم توضیح داده می‌شود، به آن چه می‌گویند؟'
$questions[202] = 'یک دختر ۴.۵ ساله همیشه حتی در فصل تابستان مجبور به پوشیدن جوراب گرم بود. در معاینه فیزیکی مشخص شد که او فشار خون بالا دارد و نبض فمورال او در مقایسه با نبض رادیال و کاروتید ضعیف بود. رادیوگرافی قفسه سینه نشان دهنده فرورفتگی قابل توجه دنده‌ها در امتداد لبه‌های پایینی آن‌ها بود. این وضعیت به دلیل چه چیزی بود؟'
$questions[203] = 'علامت فاژه در کدام بیماری دیده می‌شود؟'
$questions[204] = 'پروتوزوئای مرتبط با مگاازوفاگوس -'
$questions[205] = 'کدام یک از گیرنده های حسی زیر در اپیدرم یافت می شوند؟'
$questions[206] = 'بیمار با دیپلوپی تک چشمی به بخش اورژانس مراجعه می‌کند. معاینه با نور مایل یک هلال طلایی و معاینه با نور محوری یک خط هلالی تیره نشان می‌دهد. کدام یک از موارد زیر محتمل‌ترین تشخیص است؟'
$questions[207] = 'شایع‌ترین نوع آترزی مادرزادی چیست؟'
$questions[208] = 'هاراکیری مرگ با چیست؟'
$questions[209] = 'کدام یک از کمان های حلقی زیر در دوران جنینی منجر به تشکیل اپی گلوت می شود؟'
$questions[210] = 'مرکز تنفس توسط همه موارد زیر به جز ... دچار افسردگی می‌شود.'
$questions[211] = 'نوع شکستگی در استخوان صخره‌ای'
$questions[212] = 'کدام یک از آنزیم های زیر مرحله غیرقابل برگشت در گلیکولیز را کاتالیز نمی کند؟'
$questions[213] = 'تعداد کل دهیدروژنازهای چرخه کربس؟'
$questions[214] = 'کارسینوئید قلب به چه صورت تظاهر می‌کند؟'
$questions[215] = 'ظاهر شعاع‌مانند (sunburst) در رادیوگرافی در کدام مورد دیده می‌شود؟'
$questions[216] = 'مناسب‌ترین دارو برای القا و نگهداری داخل وریدی در جراحی روزانه کدام است؟'
$questions[217] = 'یک زن آسیایی ۲۲ ساله مهاجر با شکایت از ضعف، تب، درد بازو، کاهش اشتها و مشکلات بینایی مراجعه می‌کند. مادرش بیان می‌کند که یک هفته پیش غش کرده است. پزشک قادر به لمس نبض‌های اندام تحتانی بیمار نیست و متوجه می‌شود که نبض‌های مچ دست ضعیف هستند. میزان رسوب گلبول‌های قرمز (esr) افزایش یافته است. کدام یک از موارد زیر محتمل‌ترین تشخیص است؟'
$questions[218] = 'گانگلیون ستلاره در کدام ناحیه دیده می‌شود؟'
$questions[219] = 'پیشگیری پس از مواجهه با hiv باید در اسرع وقت پس از مواجهه و حداکثر ظرف مدت ... آغاز شود.'
$questions[220] = 'در مورد کراتوآکانتوما کدام گزینه صحیح است؟'
$questions[221] = 'صدای مردانه در زنان با کدام روش درمان می‌شود؟'
$questions[222] = 'اثربخشی دارو به چه معناست؟'
$questions[223] = 'در مورد سرطان نازوفارنکس همه موارد زیر صحیح است به جز'
$questions[224] = 'بررسی مغز استخوان یک کودک ۲ ساله با طحال بسیار بزرگ و پانسیتوپنی، سلول‌هایی با ظاهر کاغذ چروکیده در سیتوپلاسم نشان داد. تشخیص احتمالی چیست؟'
$questions[225] = 'درمان خوراکی انتخابی برای عفونت پوستی ناشی از mrsa کدام است؟'
$questions[226] = 'در مورد رینوره csf کدام گزینه صحیح است؟ up 09'
$questions[227] = 'در پایان هفته پنجم بارداری، چند عدد سومیت قابل مشاهده است؟'
$questions[228] = 'تخلیه لنفاوی دهانه رحم توسط تمام گره های لنفاوی زیر انجام می شود، به جز کدام یک؟'
$questions[229] = 'علامت prehn در کدام مورد مشاهده می‌شود؟'
$questions[230] = 'سابقه عدم علاقه به غذاهای شیرین معمولاً در کدام یک از موارد زیر وجود دارد؟'
$questions[231] = 'هموپکسین به کدام یک متصل می‌شود؟'
$questions[232] = 'نتایج مثبت کاذب در تست vdrl در کدام مورد دیده می‌شود؟'
$questions[233] = 'کدام هورمون در رفلکس ترشح شیر نقش دارد؟'
$questions[234] = 'کدام بی‌حس‌کننده موضعی به آدرنالین نیاز ندارد؟'
$questions[235] = 'یک بیمار با عملکرد کلیوی طبیعی به مدت ۲ هفته دوز نگهدارنده روزانه دیگوکسین دریافت کرده است. اگر دوز تغییر کند، غلظت پلاسمایی دیگوکسین در حالت پایدار جدید تقریباً در چه مدت زمانی حاصل می‌شود؟'
$questions[236] = 'کدام یک از موارد زیر با سه گانه ویرکوف مرتبط است؟'
$questions[237] = 'فئوکروموسیتوما با ترشح بیش از حد کدام یک از موارد زیر مشخص می‌شود؟'
$questions[238] = 'کدام گزینه در مورد کیست آدنوم موکینوز پانکراس نادرست است؟'
$questions[239] = 'جووار به دلیل بیش‌بود کدام ماده باعث ایجاد پلاگر می‌شود؟'
$questions[240] = 'ترومبوز شریانی در کدام یک از موارد زیر دیده می‌شود؟'
$questions[241] = 'یک مرد با پنوموتوراکس متوسط سمت راست بدون تنش مراجعه کرده است، یافته های فیزیکی به صورت زیر است:'
$questions[242] = 'پنومونیت حساسیتی به طور کلاسیک به عنوان کدام نوع واکنش حساسیتی توصیف می‌شود؟'
$questions[243] = 'مدیریت اولیه آنتی ژن مادرزادی عدسی چیست؟'
$questions[244] = 'ادم روی ماستوئید در کدام مورد دیده می‌شود؟'
$questions[245] = 'شایع ترین کمبود آنزیمی مسئول گالاکتوزمی کدام است؟'
$questions[246] = 'هیدروپس جنینی غیرایمنی ناشی از کدام مورد است؟'
$questions[247] = 'بیمار مبتلا به سرطان دهان، دارای غده لنفاوی هم‌طرفی به اندازه 2 سانتی‌متر، تک‌عدد و با گسترش کپسولی است. بر اساس مرحله‌بندی tnm، این مورد در کدام دسته قرار می‌گیرد؟'
$questions[248] = 'بر اساس طبقه‌بندی الیس و دیوی، شکستگی دندان ۵۱ که شامل عاج و همراه با از دست دادن حیات دندان باشد، در کدام دسته قرار می‌گیرد؟'
$questions[249] = 'یافته صحیح در مورد ویروس هپاتیت c کدام است؟'
$questions[250] = 'ناحیه کمری ستون فقرات تمام حرکات زیر را مجاز می‌کند به جز:'
$questions[251] = 'شدت نارسایی میترال ممکن است بر اساس چه موردی ارزیابی شود؟'
$questions[252] = 'حبس تا هفت سال و همچنین جریمه به عنوان مجازات برای ایراد عمدی آسیب شدید:'
$questions[253] = 'حساسترین آزمایش برای تشخیص عفونت hiv'
$questions[254] = 'کدام یک از موارد زیر در مورد سیاه زخم (آنتراکس) صحیح است؟'
$questions[255] = 'شایع‌ترین محل متاستاز استخوان تمپورال معمولاً در کدام مورد دیده می‌شود؟'
$questions[256] = 'ظهور یا بازظهور در کدام یک از ارگانیسم های زیر مشاهده شده است -'
$questions[257] = 'در اسکن تیروئید زیر، محتمل‌ترین تشخیص کدام است؟'
$questions[258] = 'اتوکلاو کردن در چه شرایطی انجام می‌شود؟'
$questions[259] = 'یک زن 28 ساله با سوختگی در ناحیه سر، گردن و جلوی تنه به شما مراجعه می‌کند. درصد سطح سوختگی را در این سناریو محاسبه کنید.'
$questions[260] = 'واسکولیت حساسیت بیش از حد در کدام یک دیده می‌شود؟'
$questions[261] = 'سه‌گانه شارکو شامل چه مواردی است؟'
$questions[262] = 'علامت "حباب دوتایی" در رادیوگرافی در کدام مورد مشاهده می‌شود؟'
$questions[263] = 'قوی‌ترین اتصال زانول در کجاست؟'
$questions[264] = 'کدام یک از موارد زیر در وابستگی به یک ماده ابتدا ایجاد می‌شود؟'
$questions[265] = 'شایع‌ترین علامت قابل مشاهده در دررفتگی مادرزادی مفصل ران در کودکان بزرگتر چیست؟'
$questions[266] = 'فلج دوطرفه عصب فرنیک توسط کدام عامل ایجاد می‌شود؟'
$questions[267] = 'حفره میان‌مغز'
$questions[268] = 'بارزترین نشانه کم خونی مگالوبلاستیک کدام است؟'
$questions[269] = 'آزمون اولانی برای تشخیص چه موردی استفاده می‌شود؟'
$questions[270] = 'نقص میدان دید دوطرفه گیجگاهی مشخصه کدام یک از موارد زیر است؟'
$questions[271] = 'تمامی درمان‌های زیر ممکن است در یک نوزاد یک ساعته با آسفیکسی شدید مورد نیاز باشد به جز-'
$questions[272] = 'کاهش 10 درجه سانتی‌گراد دما باعث کاهش میزان متابولیسم مغز به میزان'
$questions[273] = 'بیمار با دررفتگی قدامی شانه به احتمال زیاد سابقه کدام حالت را گزارش می‌کند؟'
$questions[274] = 'در بی‌حسی دندان‌های قدامی فک پایین، همه موارد زیر به‌جز کدام یک مورد استفاده قرار می‌گیرند؟'
$questions[275] = 'یک مرد 61 ساله برای بررسی اتساع شکمی بدون دلیل، تحت سی تی اسکن شکم قرار می‌گیرد. جمع‌شدگی‌های داخل صفاقی با کاهش چگالی و سپتوم‌های تقویت‌شده مشاهده می‌شود. همچنین فرورفتگی در حاشیه کبد و ضخیم شدن امنتوم دیده می‌شود. کدام یک از موارد زیر به‌عنوان علت زمینه‌ای این یافته‌ها محتمل‌تر است؟'
$questions[276] = 'زمان اولیه سفت شدن برای خمیر قالب گیری زیر چیست؟'
$questions[277] = 'ارزیابی بیمار مبتلا به پروستاتیسم شامل همه موارد زیر به جز کدام است؟'
$questions[278] = 'داروی a دارای زیست دسترسی 80٪ و حجم توزیع 10 لیتر است. دوز بارگیری مورد نیاز برای رسیدن به غلظت پلاسمایی 0.6 میلی‌گرم بر لیتر را محاسبه کنید.'
$questions[279] = 'یک بیمار پس از عمل دچار سپتیسمی شد و به صورت تجربی توسط یک پزشک مقیم جدید تحت درمان ترکیبی شیمی درمانی قرار گرفت. با این حال، هنگامی که بیمار حتی پس از ۱۰ روز درمان با آنتی‌بیوتیک پاسخ نداد، بررسی پرونده‌ها انجام شد. مشخص شد که پزشک مقیم ترکیبی از آنتی‌بیوتیک‌ها را شروع کرده بود که اثر متقابل آنتاگونیستی داشتند. کدام یک از ترکیبات زیر به احتمال زیاد تجویز شده بود؟'
$questions[280] = 'رحول با همسایه خود دعوا کرد و به او حمله کرد. هیچ آسیبی وجود نداشت، اما همچنان می‌توان او را تحت کدام بخش قانونی پیگرد قرار داد؟ aiims 12'
$questions[281] = 'لکه‌های میلیاری در کدام یک از موارد زیر دیده می‌شود؟ الف) سل ب) سارکوئیدوز ج) سیلیکوز د) پنومونی ناشی از p. carinii'
$questions[282] = 'شل‌کننده عضلانی انتخابی در بیمار با بیلیروبین سرمی 6 میلی‌گرم در دسی‌لیتر و کراتینین سرمی 4.5 میلی‌گرم در دسی‌لیتر کدام است؟'
$questions[283] = 'کدام لایه شبکیه در تماس با زجاجیه است؟'
$questions[284] = 'شایع ترین تومور بدخیم اولیه استخوان -'
$questions[285] = 'نوع دوم اسکیزوفرنی با همه موارد زیر به جز کدامیک مشخص می‌شود؟'
$questions[286] = 'نرخ خالص تولید مثل برابر با 1 به معنای نرخ حفاظت زوجین برابر با چه مقداری است؟ (تکرار)'
$questions[287] = 'کدام یک از داروهای زیر در برابر ارگانیسمی که آنزیم های غیرفعال کننده آمینوگلیکوزید تولید می‌کند، موثرتر است؟'
$questions[288] = 'وظیفه محول شده به راهنمای سلامت روستایی: مارس 2013 (ب)'
$questions[289] = 'کدام بیماری توسط آنتی بادی ضد فسفولیپاز رسیپتور ایجاد می شود؟'
$questions[290] = 'کدام یک از موارد زیر در مورد لیگامان طولی خلفی استخوانی شده نادرست است؟'
$questions[291] = 'رحم حلزونی چیست؟'
$questions[292] = 'کدام یک از موارد زیر یک بی‌حسی موضعی با پیوند استری است؟'
$questions[293] = 'کدام یک از موارد زیر باعث استریدور در نوزادان نمی‌شود؟'
$questions[294] = 'ضربه زایمان یک عامل خطر برای کدام مورد است؟'
$questions[295] = 'رنگ‌آمیزی ویژه مورد استفاده برای ترپونما پالیدوم کدام است؟'
$questions[296] = 'کدام عصب در شکستگی گردن فیبولا آسیب می‌بیند؟'
$questions[297] = 'بیماری اسگود شلاتر شامل کدام قسمت می‌شود؟'
$questions[298] = 'نتیجه مثبت کاذب تست توبرکولین نشان دهنده چیست؟'
$questions[299] = 'بیماری ولز توسط چه چیزی ایجاد می‌شود؟'
$questions[300] = 'یک پسر 4 ساله با عقب ماندگی ذهنی، خودزنی و هیپراوریسمی به احتمال زیاد دچار کمبود آنزیمی است که در کدام فرآیند نقش دارد؟'
$questions[301] = 'همه موارد زیر در سوماتیزاسیون هیستریک دیده می‌شوند به جز:'
$questions[302] = 'تومور سلول های آسینار بیشتر در کدام محل دیده می شود؟'
$questions[303] = 'یک زن 28 ساله با شکایت از کاهش خواب و رفتار آشفته به مدت 9 ماه مراجعه می‌کند. او احساس می‌کند که یک دوربین پشت سرش نصب شده که همیشه او را دنبال می‌کند. کدام یک از موارد زیر احتمالاً تشخیص صحیح است؟'
$questions[304] = 'چه ناهنجاری در موهای یک کودک مبتلا به سوءتغذیه شدید مشاهده می‌شود؟'
$questions[305] = 'بر اساس طبقه‌بندی who برای سرکوب سیستم ایمنی، سطح cd4 در کودکان مبتلا به hiv/aids در مرحله شدید در گروه سنی 36 تا 59 ماه چقدر است؟'
$questions[306] = 'انتقال کلرید به سرعت رخ می‌دهد و اساساً در چه مدت زمانی کامل می‌شود؟'
$questions[307] = 'خونرسانی به کولون سیگموئید توسط کدام شریان انجام می‌شود؟  
الف) شریان کولیک میانی  
ب) شریان حاشیه‌ای  
ج) شریان کولیک چپ  
د) شریان سیگموئید'
$questions[308] = 'شایع‌ترین علت لوکوریا در نوزادان چیست؟'
$questions[309] = 'تمامی ساختارهایی که از روزنه مشخص شده در نمودار عبور می کنند به جز کدام مورد هستند؟'
$questions[310] = 'همه موارد زیر در ترک نیکوتین مشاهده می‌شوند به جز؟'
$questions[311] = 'دوز ویتامین a برای یک کودک 18 ماهه مبتلا به کراتومالاسی، با وزن 10 کیلوگرم چقدر است؟'
$questions[312] = 'محل‌های گیرنده نیکوتینی شامل همه موارد زیر به جز کدام است؟'
$questions[313] = 'شایع‌ترین علت بدخیمی استخوانی -'
$questions[314] = 'کدام لایه پوست به عنوان سدی در برابر از دست دادن آب عمل می‌کند؟'
$questions[315] = 'آنژیوگرافی ایندوسیانین گرین (icg آنژیوگرافی) بیشترین کاربرد را در تشخیص کدام مورد دارد؟'
$questions[316] = 'مخرج برای محاسبه نرخ مرگ و میر نسبی از یک بیماری خاص چیست؟'
$questions[317] = 'در انسان، عرض کورتکس مو معمولاً:'
$questions[318] = 'شیردهی توسط کدام مورد مهار می‌شود؟'
$questions[319] = 'کلاس چهارم طبقه‌بندی کندی چیست؟'
$questions[320] = 'همه موارد زیر از تظاهرات کلاسیک ناهنجاری‌های اتصال جمجمه-مهره‌ای هستند به جز'
$questions[321] = 'رفتار شناختی با کدام مورد سر و کار دارد؟  
الف) انگیزه ناخودآگاه  
ب) فرضیات ناسازگار  
ج) تعارض شبه پویا  
د) آگاهی هیجانی از افکار خودکار'
$questions[322] = 'گزینه‌های درست درباره نوكاردیا به جز کدام است؟'
$questions[323] = 'شایع‌ترین محل زخم کورلینگ کدام است؟'
$questions[324] = 'بیوپسی کبد از طریق خط میانی زیر بغل هشتم انجام می‌شود تا از چه چیزی اجتناب شود؟'
$questions[325] = 'یک نوزاد تازه متولد شده در مدت 48 ساعت مکونیوم دفع نکرده است. روش تشخیصی انتخابی چیست؟'
$questions[326] = 'کدام یک از شریان‌های زیر شاخه مستقیم شریان گاسترودوئودنال است؟'
$questions[327] = 'کدام یک از عضلات زیر در اثر آسیب به عصب مدیان در مچ دست فلج نمی‌شود؟'
$questions[328] = 'انتقال عمودی آنتی‌بادی می‌تواند در همه موارد زیر رخ دهد، به جز -'
$questions[329] = 'فردی با نارسایی میترال و فیبریلاسیون دهلیزی دچار سنکوپ شده است. در معاینه، ضربان قلب فرد 55 است. محتمل‌ترین علت چیست؟'
$questions[330] = 'در فلج دورکننده تارهای صوتی دوطرفه (b/l)، کدام یک از موارد زیر انجام نمی‌شود؟'
$questions[331] = 'آرتریت تام اسمیت به چه علتی است؟'
$questions[332] = 'کدام یک از موارد زیر در مورد وبا صحیح است؟'
$questions[333] = 'شایع ترین ناهنجاری رحمی کدام است؟ مارس 2007'
$questions[334] = 'تست ''نشاسته ید'' برای تشخیص کدام مورد مفید است؟'
$questions[335] = 'از فر هوای گرم برای استریل کردن چه چیزی استفاده می‌شود؟'
$questions[336] = 'در مورد عفونت hiv، کدام گزینه نادرست است؟'
$questions[337] = 'فشارهای دمی پایدار بالا در طول تهویه با فشار مثبت، خطر کدام مورد را افزایش می‌دهد؟'
$questions[338] = 'مثال آپوپتوز چیست؟'
$questions[339] = 'تریاد تروتر شامل همه موارد زیر به جز کدام است؟'
$questions[340] = 'سریع‌ترین داروی سیکلوپلژیک کدام است؟'
$questions[341] = 'عامل خطر در کولانژیوکارسینوما کدام است؟'
$questions[342] = 'یک پسر 9 ساله با سابقه کاهش حجم ادرار، ادرار به رنگ کولا و تورم صورت و دست ها به مدت 2 روز مراجعه می کند. او مبتلا به فشار خون بالا، صورت پف کرده و ادم فرورونده اندام تحتانی است. سابقه ضایعات پوستی 4 هفته قبل را دارد. تشخیص گلومرولونفریت پس از استرپتوکوکی داده می شود. تیتر aslo احتمالاً چگونه است؟'
$questions[343] = 'بهترین روش تصویربرداری در کارسینوم برونکوژنیک کدام است؟'
$questions[344] = 'در بیماری ادرار شربت افرا، همه اسیدهای آمینه زیر در ادرار دفع می‌شوند، به جز:'
$questions[345] = 'الگوی توارث کانتراکچر دوپویترن به چه صورت است؟'
$questions[346] = 'بزرگ شدن غدد لنفاوی اینگوئینال در کدام مورد مشاهده می‌شود؟'
$questions[347] = 'اعلامیه ژنو مربوط به کدام مورد است؟'
$questions[348] = 'در مورد پروب نابر کدام گزینه نادرست است؟'
$questions[349] = 'همه موارد زیر گشادکننده عروق هستند به جز:'
$questions[350] = 'تب خونریزی‌دهنده ناشی از ویروس‌ها شامل همه موارد زیر است به جز:'
$questions[351] = 'تمامی تغییرات عروقی زیر در التهاب حاد مشاهده می‌شوند، به جز'
$questions[352] = 'اختلال بافت همبند مختلط با کدام یک از موارد زیر مرتبط نیست؟'
$questions[353] = 'آنوریسم محیطی -'
$questions[354] = 'راسبوریکاز یک داروی جدیدتر مورد استفاده در نقرس است. این دارو با کدام مکانیسم عمل می‌کند؟'
$questions[355] = 'پاسخ ادیومتری امپدانس در اوتیت مدیا سروز:'
$questions[356] = 'سلول های ستاره ای ون کوفر در سینوزوئیدهای کدام یک از اندام های زیر دیده می شوند؟'
$questions[357] = 'محصول نهایی تغییر نیافته چرخه tca کدام است؟'
$questions[358] = 'اتصال اولیه بین چرخه اسید سیتریک و چرخه اوره با کدام است؟'
$questions[359] = 'آتروفی قهوه ای ناشی از تجمع کدام ماده است؟'
$questions[360] = 'تمام مراکز استخوان‌سازی اولیه در چه سن جنینی ظاهر می‌شوند؟ wb 11'
$questions[361] = '"ارزیابی گنادوتروپین" برای ارزیابی، مناسب‌ترین روز از یک چرخه قاعدگی 28 روزه طبیعی برای زنی با دوره قاعدگی 5 روزه را انتخاب کنید. (1 روز را انتخاب کنید)'
$questions[362] = 'علل آلکالوز متابولیک شامل همه موارد زیر به جز کدام است؟'
$questions[363] = 'اولین علامت سرطان فرج چیست؟'
$questions[364] = 'غدد عرق از کدام قسمت عصب دهی کولینرژیک دریافت می‌کنند؟'
$questions[365] = 'رویش اولین دندان آسیای بزرگ دائمی در چه سنی اتفاق می‌افتد؟'
$questions[366] = 'تولید hmg co a در کبد توسط کدام مورد مهار می‌شود؟'
$questions[367] = '"اجسام کال-اکزنر" در کدام مورد دیده می‌شوند؟'
$questions[368] = 'متاستاز ریه به ریه در کدام مورد مشاهده می‌شود؟'
$questions[369] = 'فلج کلامپکه به دلیل آسیب به کدام یک از اعصاب زیر رخ می‌دهد؟'
$questions[370] = 'کدام یک از داروهای ضد سرطان زیر باعث نوروپاتی محیطی به عنوان یک عارضه جانبی می‌شود؟'
$questions[371] = 'بر اساس نام‌گذاری کویناد، کدام یک از بخش‌های زیر کبد دارای عروق‌خونی مستقل است؟'
$questions[372] = 'تکمیل تاج دندان شیری مولر اول فک پایین در چه زمانی اتفاق می‌افتد؟'
$questions[373] = 'اجسام پساموما در کدام یک دیده می‌شوند؟'
$questions[374] = 'یک مرد 58 ساله برای ارزیابی علائم درد قفسه سینه تحت کاتتریزاسیون قلبی قرار می‌گیرد. او نگران خطرات این روش است و به عنوان بخشی از رضایت آگاهانه، شما او را در مورد خطرات و مزایای این روش راهنمایی می‌کنید. کدام یک از جنبه‌های آنژیوگرافی زیر صحیح است؟'
$questions[375] = 'کدام دارو در متابولیسم پیریدوکسین اختلال ایجاد می‌کند؟'
$questions[376] = 'اجسام وایبل-پالاد در کدام سلول‌ها وجود دارند؟'
$questions[377] = 'کدام یک از اندام های زیر احتمالاً دچار نکروز انعقادی نمی شوند؟'
$questions[378] = 'جنسیت جنین از روی اندام تناسلی خارجی در چه زمانی به وضوح قابل تشخیص می‌شود؟'
$questions[379] = 'خونریزی واکنشی در چه زمانی رخ می‌دهد؟'
$questions[380] = 'یک مورد تروما غیر نافذ به اورژانس آورده شده است، در حالت شوک است؛ به کریستالوئیدهای داخل وریدی پاسخ نمی‌دهد؛ مرحله بعدی در مدیریت او چیست؟'
$questions[381] = 'سیتوپاتولوژی به چه چیزی می‌پردازد؟'
$questions[382] = 'رشته‌های حسی جوانه‌های چشایی در زبان و کام نرم از طریق کدام عصب منتقل می‌شوند؟'
$questions[383] = 'یافته های هیستولوژیک رد حاد پیوند کلیه کدامند؟'
$questions[384] = 'سنی که کودک می‌تواند جنسیت را تشخیص دهد؟'
$questions[385] = 'دو سوم قدامی زبان از کدام بخش تشکیل میشود؟'
$questions[386] = 'افزایش تب باعث افزایش دفع آب به میزان ______ میلی‌لیتر در روز به ازای هر درجه سانتی‌گراد می‌شود.'
$questions[387] = 'کدام سلول در لنفوم هوچکین دیده نمی‌شود؟'
$questions[388] = 'سندرم کودک ناجور در کدام دسته از اختلالات زیر طبقه بندی می شود؟'
$questions[389] = 'کمبود کدام آنزیم منجر به سیترولینمی نوع 1 می‌شود؟'
$questions[390] = 'یک زن ۱۹ ساله در طول قاعدگی دچار خونریزی از بینی می‌شود. به احتمال زیاد چه بیماری دارد؟'
$questions[391] = 'بیماری درکوم شایع‌ترین در کدام قسمت است؟'
$questions[392] = 'در دوران بارداری، فیبروم ممکن است تمام عوارض زیر را داشته باشد به جز: مارس 2009'
$questions[393] = 'کدام بخش از پوشش داخلی رحم که در طول قاعدگی ریزش نمی‌کند، است؟'
$questions[394] = 'آبسه زیر بغل به چه روشی به طور ایمن تخلیه می‌شود؟'
$questions[395] = 'پروتئین های c-reactive جزو کدام دسته هستند؟'
$questions[396] = 'بهترین مدیریت برای بیمار با وضعیت همودینامیک پایدار که نوار قلب او کمپلکس qrs پهن و تاکیکاردی آنتی درومیک نشان می‌دهد چیست؟'
$questions[397] = 'حفظ آبزدایی قرنیه توسط کدام بخش انجام می‌شود؟'
$questions[398] = 'کدام یک از یافته‌های زیر در معاینه یک کودک با شکاف کام درجه سه مشاهده می‌شود؟'
$questions[399] = 'یک زن ۲۷ ساله نولیپار از ۴ ماه پیش از منوراژی شدید و درد زیر شکم شکایت دارد. در معاینه، رحم به اندازه ۹ هفته همراه با فیبروئید فوندال مشاهده می‌شود. درمان انتخابی چیست؟'
$questions[400] = '"اجسام کریولا" در خلط نشانه پاتوگنومونیک کدام بیماری است؟'
$questions[401] = 'نیروی نگهدارنده اورژانسی برای دنچر کامل فک بالا چیست؟'
$questions[402] = 'تمام موارد زیر از نگهدارنده‌های قوی رحم هستند به جز:'
$questions[403] = 'موارد زیر همگی از اختلالات میلوپرولیفراتیو هستند به جز-'
$questions[404] = 'خانمی با ناباروری و انسداد دوطرفه لوله‌های رحمی در ناحیه کورنوآ، بهترین روش مدیریت چیست؟'
$questions[405] = 'کارنیتین از چه ترکیباتی ساخته شده است؟'
$questions[406] = 'ریتم eeg ثبت شده از سطح پوست سر در طول خواب rem کدام است؟'
$questions[407] = 'دوسوگرایی (ambivalence) بیشتر در کدام یک از موارد زیر مشاهده می‌شود؟'
$questions[408] = 'متابولیسم پرولین در کجا انجام می‌شود؟'
$questions[409] = 'یک پسر ۳ ساله با شروع ناگهانی کاهش دید در چشم چپ مراجعه می‌کند. معاینه چشم‌پزشکی، استرابیسم و رفلکس سفید مردمک را نشان می‌دهد. یافته‌های فوندوس در زیر نمایش داده شده است. تشخیص شما چیست؟'
$questions[410] = 'یک مرد ۵۰ ساله با راه رفتن "تلوتلو خوران" و درد "برق‌آسا" در دست‌ها و پاها مراجعه می‌کند. سابقه پزشکی او شامل آنوریسم آئورت و نارسایی آئورت است. معاینه عصبی اختلال در حس لرزش، لمس و درد در اندام تحتانی را نشان می‌دهد. بیمار متعاقباً به دلیل ذات‌الریه فوت می‌کند. کالبدشکافی، التهاب انسدادآور شریان‌های مننژ و آتروفی ستون‌های خلفی نخاع را آشکار می‌سازد. تشخیص مناسب چیست؟'
$questions[411] = 'وزن تقریبی غده تیموس در زمان بلوغ چقدر است؟'
$questions[412] = 'اجسام لوی حاوی چه چیزی هستند؟'
$questions[413] = 'یک دختر ۱۲ ساله با قد کوتاه، گردن پرده‌دار و فاصله زیاد نوک پستان‌ها مراجعه می‌کند که نشان‌دهنده یک ناهنجاری کروموزومی است و با کاریوتایپ تأیید می‌شود. کدام ناهنجاری قلبی-عروقی به احتمال زیاد در این کودک وجود دارد؟'
$questions[414] = 'سوراخ شدن متعدد پرده صماخ مشخصه کدام یک از موارد زیر است؟'
$questions[415] = 'یک پسر ۱۵ ساله با سابقه یک روزه خونریزی لثه، خونریزی زیر ملتحمه و راش پورپوریک مراجعه کرده است. بررسی‌ها نتایج زیر را نشان داد:
hb-6.4 گرم در دسی‌لیتر
tlc-26,500 در میلی‌متر مکعب؛ پلاکت - 35,000 در میلی‌متر مکعب
زمان پروترومبین-20 ثانیه با کنترل 13 ثانیه
زمان ترومبوپلاستین جزئی-50 ثانیه و
فیبرینوژن 10 میلی‌گرم در دسی‌لیتر
اسمیر محیطی نشان‌دهنده لوسمی میلوبلاستیک حاد بود. کدام یک از گزینه‌های زیر محتمل‌ترین است؟'
$questions[416] = 'کدام یک از داروهای زیر حتی کمترین اثر آگونیستی را ندارد؟'
$questions[417] = 'سرطان کیسه بیضه به طور سنتی در کدام شغل دیده می‌شود؟'
$questions[418] = 'واکنش دو فازی در فشار خون با تجویز کدام یک از موارد زیر مشاهده می‌شود؟'
$questions[419] = 'شایع‌ترین فرم سل خارج ریوی کدام است؟'
$questions[420] = 'کدام یک از مهارکننده های پمپ پروتون زیر دارای فعالیت مهارکنندگی آنزیمی است؟'
$questions[421] = 'میزان اکسیژن در مخلوط بیهوشی چقدر است؟'
$questions[422] = 'سندرم شبه تومور در کدام یک از موارد زیر دیده می‌شود؟'
$questions[423] = 'اولین سد شیمیایی که میکروارگانیسم در محل‌های شایع در معرض مواجهه با آن قرار می‌گیرد چیست؟'
$questions[424] = 'حسادت بیمارگونه با کدام مورد مرتبط است؟'
$questions[425] = 'برای پمپ سدیم-پتاسیم، نسبت کوپلینگ چیست؟'
$questions[426] = 'کدام یک از داروهای ضد رتروویروسی زیر باعث دیس لیپیدمی نمی‌شود؟'
$questions[427] = 'بیمار از دیدن هاله‌های رنگی در عصر و تاری دید در چند روز اخیر با فشار داخل چشمی طبیعی شکایت دارد:'
$questions[428] = 'چه کسی فرض کرد که فلورایددار کردن آب با کاهش پوسیدگی دندان مرتبط است؟'
$questions[429] = 'چرا به آن cu t 200 می‌گویند؟'
$questions[430] = 'نقش توسیلیزوماب در بیمار مبتلا به سندرم دیسترس تنفسی حاد (ards) ناشی از کووید-19 چیست؟'
$questions[431] = 'پروتئین‌های بنس جونز چه هستند؟'
$questions[432] = 'میزان حمله ثانویه آبله مرغان چقدر است؟'
$questions[433] = 'عبارت صحیح در مورد سیلیکوزیس -'
$questions[434] = 'از bcg برای چه مواردی استفاده می‌شود؟'
$questions[435] = 'کدام یک از گزینه‌های زیر مربوط به دسته‌بندی 3 در طبقه‌بندی ماستریخت برای اهدا پس از مرگ قلبی است؟'
$questions[436] = 'مدفوع دودی رنگ در کدام مسمومیت دیده می‌شود؟'
$questions[437] = 'کم خونی مگالوبلاستیک در کدام مورد مشاهده می‌شود؟'
$questions[438] = 'اسید آمینه در انتهای کاهنده گلوتاتیون کدام است؟'
$questions[439] = 'نوروترانسمیتر اصلی در آوران‌های هسته منفرد (nucleus tractus solitarius) برای تنظیم سیستم قلبی عروقی کدام است؟'
$questions[440] = 'کدام یک از موارد زیر ژن تعیین کننده تخمدان است؟'
$questions[441] = 'اولین عضله درگیر در افتالموپاتی تیروئیدی کدام است؟'
$questions[442] = 'تب اسکارلت توسط کدام عامل ایجاد می‌شود؟'
$questions[443] = 'معاینه ستون فقرات در بیمار دچار تروماهای متعدد چگونه انجام می‌شود؟'
$questions[444] = 'کدام یک از موارد زیر در مورد آنژیوفیبروم نوجوانان صحیح نیست؟'
$questions[445] = 'در مورد گره‌های رانویه کدام گزینه صحیح است؟'
$questions[446] = 'افزایش تعداد wbc در تمام موارد زیر دیده می‌شود به جز؟'
$questions[447] = 'پس از رابطه جنسی، فرد در بیضه چپ خود درد ایجاد می‌کند که با بالا آوردن کیسه بیضه تسکین نمی‌یابد. تشخیص چیست؟'
$questions[448] = 'طول سیگموئیدوسکوپ انعطاف پذیر ......... سانتی متر است:'
$questions[449] = 'عامل ایجاد شانکروئید کدام است؟'
$questions[450] = 'کاهش وزن قابل توجه عبارت است از -'
$questions[451] = 'علائم چوستک و تروسو در کدام مسمومیت دیده می‌شود؟'
$questions[452] = 'داروهای سمپاتومیمتیک در درمان تمامی موارد زیر مفید هستند به جز -'
$questions[453] = 'تفاوت بین زخم بریده شده و پارگی ناشی از حرارت یا پارگی ناشی از حرارت چیست؟'
$questions[454] = 'درمان انتخابی برای کارسینومای کوچک پنیس در ناحیه پوست ختنه‌گاه چیست؟'
$questions[455] = 'شریان مهره‌ای از سوراخ زائده‌های عرضی کدام مهره‌ها عبور می‌کند؟'
$questions[456] = 'مأموریت ایندرادانوش برای چیست؟'
$questions[457] = 'بیماری منتقله از آب کدام است؟'
$questions[458] = 'یک زن مسن با سابقه بیماری الکلی دچار زردی و ادم شدید (آناسارکا) می‌شود. کدام یک از موارد زیر محتمل‌ترین پاتوفیزیولوژی برای ادم پایدار او است؟'
$questions[459] = 'شاخص درمانی معیاری است برای'
$questions[460] = 'کدام یک از روش‌های زیر بهترین روش برای تعیین عمق تهاجم در کارسینومای مری است؟'
$questions[461] = 'رگ‌های موجود در بند ناف کدامند؟'
$questions[462] = 'برونکوگرافی ممکن است خطرناک باشد اگر بیمار دارای:'
$questions[463] = 'یک زن 40 ساله از سوزش سر دل در ناحیه اپیگاستر و رترواسترنال شکایت دارد. او همچنین علائم رگورژیتاسیون را تجربه می‌کند. آندوسکوپی اریتم مری را نشان می‌دهد که با ازوفاژیت رفلاکسی سازگار است. بیمار اقدامات محافظه‌کارانه از جمله مصرف مهارکننده‌های پمپ پروتون (ppi) را امتحان کرده اما بهبودی در علائم مشاهده نشده است. کدام یک از گزینه‌های زیر صحیح است؟'
$questions[464] = 'یک کودک 2 ساله با فرو رفتگی بین دنده‌ای و افزایش سیانوز با سابقه آسپیراسیون جسم خارجی آورده شده است. کدام اقدام ممکن است در این وضعیت نجات‌بخش باشد؟'
$questions[465] = 'تمام موارد زیر در درمان انورزیس استفاده می‌شوند به جز'
$questions[466] = 'اثر انگشت dna بر اساس وجود چه چیزی در dna است؟'
$questions[467] = 'درمان انتخابی برای انسفالوپاتی هاشیموتو -'
$questions[468] = 'نوع غضروف مشاهده شده در مفصل تمپورومندیبولار چیست؟'
$questions[469] = 'ipc 304b مربوط به چیست؟'
$questions[470] = 'اولین فاکتور اکلوژن کدام است؟'
$questions[471] = 'هپارین در بیماران مبتلا به بیماری‌های زیر منع مصرف دارد به جز'
$questions[472] = 'pseudouridine در کدام یک یافت می‌شود؟'
$questions[473] = 'عبارت صحیح در مورد کلوئید شامل همه موارد زیر به جز کدام است؟'
$questions[474] = 'به عنوان یک عارضه جانبی، سندرم متابولیک بیشتر با کدام گروه از داروهای زیر مرتبط است؟'
$questions[475] = 'این یک واقعیت رایج است که کروماتین جنسی در اجساد تجزیه شده به خوبی قابل تشخیص نیست. تا چه مدت پس از مرگ می‌توان کروموزوم y را در پالپ دندان با استفاده از رنگ‌های فلورسنت نشان داد؟'
$questions[476] = 'یک بیمار دیابتی با سلولیت اوربیتال و سینوزیت ماگزیلاری در میکروسکوپ، هیف های شفاف، باریک، سپتاته و منشعب با تهاجم به عروق خونی را نشان می‌دهد. قارچ مسبب کدام است؟'
$questions[477] = 'رشد در یک محیط جامد مصنوعی بدون سلول برای موارد زیر امکان‌پذیر است به جز -'
$questions[478] = 'کدام یک از مسدودکننده‌های h1 کمترین اثر آرام‌بخشی را دارد؟'
$questions[479] = 'در مورد میکوزیس فونگوئید کدام گزینه نادرست است؟'
$questions[480] = 'تمرین تقویت عضله'
$questions[481] = 'واکسیناسیون پس از مواجهه برای سرخک باید در چه زمانی انجام شود؟ سپتامبر 2009'
$questions[482] = 'پیشگیری شیمیایی در بیماری وبا'
$questions[483] = 'همه موارد زیر در مورد ویروس‌های تنفسی صحیح هستند به جز'
$questions[484] = 'گاستروپاتی هیپرتروفیک در کدام مورد دیده می‌شود؟'
$questions[485] = 'یک کودک مبتلا به اسهال دارای تنفس عمیق و سریع است. تشخیص چیست؟'
$questions[486] = 'سندرم مرتبط با افزایش خطر ابتلا به لوسمی کدام است؟'
$questions[487] = 'واکنش انتقال خون mc چیست؟'
$questions[488] = 'کدام یک از موارد زیر جزئی از سندرم hellp نیست؟'
$questions[489] = 'پس از پارگی شریان مننژ میانی، خونریزی در کدام ناحیه رخ می‌دهد؟'
$questions[490] = 'یک زن ۲۳ ساله با تشخیص پرولاپس دریچه میترال در اکوکاردیوگرافی برای ارزیابی سوفل سیستولیک مراجعه کرده است. کدام یک از موارد زیر بیشترین ویژگی پرولاپس دریچه میترال را دارد؟'
$questions[491] = 'استئوسارکوم از کدام نوع سلول‌ها منشأ می‌گیرد؟'
$questions[492] = 'هایپرکلسوریا در کودکان زمانی تشخیص داده می‌شود که دفع روزانه کلسیم در ادرار بیشتر از چه مقدار باشد؟'
$questions[493] = 'کدام گزینه در مورد تنگی پیلور هیپرتروفیک نادرست است؟'
$questions[494] = 'در مورد آرام‌بخشی هوشیار همه موارد زیر صحیح است به جز'
$questions[495] = 'ipc 312 و 315 مرتبط با چه مواردی هستند؟'
$questions[496] = 'همه توسط آنسا سرویکالیس عصب‌دهی می‌شوند به جز -'
$questions[497] = 'در رادیوگرافی قفسه سینه، الگوی شکل ۸ دیده می‌شود در کدام یک از موارد زیر؟'
$questions[498] = 'اجسام آشف در کدام یک از موارد زیر دیده می‌شوند؟'
$questions[499] = 'هنگام گرفتن رادیوگرافی بایت وینگ، زاویه‌دهی باید چقدر باشد تا از همپوشانی کاسپ‌ها روی سطح اکلوزال جلوگیری شود؟'
foreach ($r in $questions.Keys) {
  $ws.Cells.Item($r, 1).Value = $questions[$r]
}
